$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Regenerated s_val data (filtered save games) for sandlin_nick 2022
$data = @(
    @{Row=2; B=3.286832544864788;  C=1.655778082260271;  D=0.7527432677738641; E=0.4942365360607697; G=6.189590430959694},
    @{Row=3; B=3.286832544864788;  C=1.655778082260271;  D=3.537761648806719;  E=0.4942365360607697; G=8.974608811992548},
    @{Row=4; B=3.286832544864788;  C=1.655778082260271;  D=3.537761648806719;  E=0.4942365360607697; G=8.974608811992548},
    @{Row=5; B=1.455362044514542;  C=1.655778082260271;  D=0.1494219747398047; E=0.4942365360607697; G=3.754798637575387},
    @{Row=6; B=1.455362044514542;  C=1.655778082260271;  D=0.1494219747398047; E=10.19245300693656;  G=13.45301510845117},
    @{Row=7; B=3.286832544864788;  C=1.655778082260271;  D=0.7527432677738641; E=0.4942365360607697; G=6.189590430959694},
    @{Row=8; B=1.455362044514542;  C=1.655778082260271;  D=0.1494219747398047; E=10.19245300693656;  G=13.45301510845117},
    @{Row=9; B=3.286832544864788;  C=1.655778082260271;  D=0.1494219747398047; E=0.4942365360607697; G=5.586269137925634}
)

foreach ($row in $data) {
    $ws.Cells.Item($row.Row, 2).Value = $row.B
    $ws.Cells.Item($row.Row, 3).Value = $row.C
    $ws.Cells.Item($row.Row, 4).Value = $row.D
    $ws.Cells.Item($row.Row, 5).Value = $row.E
    $ws.Cells.Item($row.Row, 7).Value = $row.G
}
